# Apply updated "想去人数" (F column) and "最低票价" (G column) values
# to the "展览" and "全部类型" worksheets (they hold identical data).

$wb = $excel.ActiveWorkbook

# Map row number -> new F value (and optionally new G value)
$updates = @(
    @{ Row = 7;  F = 2082 },
    @{ Row = 10; F = 4553 },
    @{ Row = 15; F = 137 },
    @{ Row = 17; F = 20 },
    @{ Row = 18; F = 88 },
    @{ Row = 19; F = 3434; G = 70 },
    @{ Row = 21; F = 548 },
    @{ Row = 31; F = 679 },
    @{ Row = 32; F = 2078 },
    @{ Row = 33; F = 394 }
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates) {
        $ws.Range("F" + $u.Row).Value = $u.F
        if ($u.ContainsKey("G")) {
            $ws.Range("G" + $u.Row).Value = $u.G
        }
    }
}
